# Weekly update: insert a new price record as row 39, pushing existing
# rows 39-75 down to 40-76 (the oldest record, previously row 75, now
# becomes the new row 76).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 39..75 down by one to make room for the new weekly record.
$ws.Rows(39).Insert()

# Populate the newly inserted row 39 with this week's data.
$ws.Range("A39").Value = 2
$ws.Range("B39").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C39").Value = "Coquimbo"
$ws.Range("D39").Value = 45049
$ws.Range("E39").Value = 4
$ws.Range("F39").Value = 100112032
$ws.Range("G39").Value = "Zapallo italiano"
$ws.Range("H39").Value = "Sin especificar"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 560
$ws.Range("K39").Value = 9000
$ws.Range("L39").Value = 10000
$ws.Range("M39").Value = 9500
$ws.Range("N39").Value = "`$/caja 60 unidades"
$ws.Range("O39").Value = "Provincia de Limarí"
$ws.Range("P39").Value = 158
$ws.Range("Q39").Value = 60
$ws.Range("R39").Value = "Hortaliza"
